# Re-run of the NATMI TPM-based ligand/receptor/edge-weight computation for
# the Cadm3-Cadm3 sheet (ECs / FAPs / MuSCs cluster pairs). Ligand columns
# (G:J) depend on the sending cluster, receptor columns (M:P) depend on the
# target cluster, and edge columns (Q:T) are derived from both - this
# updates all three blocks on every data row (2-10) to the refreshed TPM
# values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.328826
$ws.Range("H2").Value = 3.986478
$ws.Range("I2").Value = 0.03246115949735631
$ws.Range("J2").Value = 0.03246115949735631
$ws.Range("M2").Value = 1.328826
$ws.Range("N2").Value = 3.986478
$ws.Range("O2").Value = 0.03246115949735631
$ws.Range("P2").Value = 0.03246115949735631
$ws.Range("Q2").Value = 1.765778538276
$ws.Range("R2").Value = 15.892006844484
$ws.Range("S2").Value = 0.001053726875912805
$ws.Range("T2").Value = 0.001053726875912805
$ws.Range("G3").Value = 1.328826
$ws.Range("H3").Value = 3.986478
$ws.Range("I3").Value = 0.03246115949735631
$ws.Range("J3").Value = 0.03246115949735631
$ws.Range("O3").Value = 0.9174921121349238
$ws.Range("P3").Value = 0.9174921121349238
$ws.Range("Q3").Value = 49.90850313826
$ws.Range("R3").Value = 449.17652824434
$ws.Range("S3").Value = 0.02978285778957808
$ws.Range("T3").Value = 0.02978285778957808
$ws.Range("G4").Value = 1.328826
$ws.Range("H4").Value = 3.986478
$ws.Range("I4").Value = 0.03246115949735631
$ws.Range("J4").Value = 0.03246115949735631
$ws.Range("M4").Value = 2.048706666666666
$ws.Range("N4").Value = 6.14612
$ws.Range("O4").Value = 0.05004672836771996
$ws.Range("P4").Value = 0.05004672836771996
$ws.Range("Q4").Value = 2.72237468504
$ws.Range("R4").Value = 24.50137216536
$ws.Range("S4").Value = 0.001624574831865424
$ws.Range("T4").Value = 0.001624574831865424
$ws.Range("I5").Value = 0.9174921121349238
$ws.Range("J5").Value = 0.9174921121349238
$ws.Range("M5").Value = 1.328826
$ws.Range("N5").Value = 3.986478
$ws.Range("O5").Value = 0.03246115949735631
$ws.Range("P5").Value = 0.03246115949735631
$ws.Range("Q5").Value = 49.90850313826
$ws.Range("R5").Value = 449.17652824434
$ws.Range("S5").Value = 0.02978285778957808
$ws.Range("T5").Value = 0.02978285778957808
$ws.Range("I6").Value = 0.9174921121349238
$ws.Range("J6").Value = 0.9174921121349238
$ws.Range("O6").Value = 0.9174921121349238
$ws.Range("P6").Value = 0.9174921121349238
$ws.Range("S6").Value = 0.8417917758298036
$ws.Range("T6").Value = 0.8417917758298036
$ws.Range("I7").Value = 0.9174921121349238
$ws.Range("J7").Value = 0.9174921121349238
$ws.Range("M7").Value = 2.048706666666666
$ws.Range("N7").Value = 6.14612
$ws.Range("O7").Value = 0.05004672836771996
$ws.Range("P7").Value = 0.05004672836771996
$ws.Range("Q7").Value = 76.94602837595555
$ws.Range("R7").Value = 692.5142553835999
$ws.Range("S7").Value = 0.04591747851554219
$ws.Range("T7").Value = 0.04591747851554219
$ws.Range("G8").Value = 2.048706666666666
$ws.Range("H8").Value = 6.14612
$ws.Range("I8").Value = 0.05004672836771996
$ws.Range("J8").Value = 0.05004672836771996
$ws.Range("M8").Value = 1.328826
$ws.Range("N8").Value = 3.986478
$ws.Range("O8").Value = 0.03246115949735631
$ws.Range("P8").Value = 0.03246115949735631
$ws.Range("Q8").Value = 2.72237468504
$ws.Range("R8").Value = 24.50137216536
$ws.Range("S8").Value = 0.001624574831865424
$ws.Range("T8").Value = 0.001624574831865424
$ws.Range("G9").Value = 2.048706666666666
$ws.Range("H9").Value = 6.14612
$ws.Range("I9").Value = 0.05004672836771996
$ws.Range("J9").Value = 0.05004672836771996
$ws.Range("O9").Value = 0.9174921121349238
$ws.Range("P9").Value = 0.9174921121349238
$ws.Range("Q9").Value = 76.94602837595555
$ws.Range("R9").Value = 692.5142553835999
$ws.Range("S9").Value = 0.04591747851554219
$ws.Range("T9").Value = 0.04591747851554219
$ws.Range("G10").Value = 2.048706666666666
$ws.Range("H10").Value = 6.14612
$ws.Range("I10").Value = 0.05004672836771996
$ws.Range("J10").Value = 0.05004672836771996
$ws.Range("M10").Value = 2.048706666666666
$ws.Range("N10").Value = 6.14612
$ws.Range("O10").Value = 0.05004672836771996
$ws.Range("P10").Value = 0.05004672836771996
$ws.Range("Q10").Value = 4.197199006044444
$ws.Range("R10").Value = 37.7747910544
$ws.Range("S10").Value = 0.002504675020312346
$ws.Range("T10").Value = 0.002504675020312346
